# Funções calcular_funcao, calcular_funcoes, chamar_calculo_das_funcoes
# Primeira versão das funções calcular_funcao, calcular_funcoes, chamar_calculo_das_funcoes
#
# Adds two new worksheets ("Funcoes_Inputs" and "Funcoes_Outputs") right after
# "Parametros" and before "Distribuições", fills them with the function
# input/output catalogue for "calcular_eventos", removes the now obsolete
# trailing blank row from "Parametros" and refreshes its AutoFilter range.

$wb = $excel.ActiveWorkbook

$parametros = $wb.Worksheets.Item("Parametros")

# --- 1. Insert the two new sheets right after Parametros, in order ---------
$wsIn = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $parametros)
$wsIn.Name = "Funcoes_Inputs"

$wsOut = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsIn)
$wsOut.Name = "Funcoes_Outputs"

# --- 2. Headers (columns A/B) for both sheets -------------------------------
$wsIn.Range("A1").Value = "Funcao"
$wsIn.Range("B1").Value = "Input"

$wsOut.Range("A1").Value = "Funcao"
$wsOut.Range("B1").Value = "Output"

# --- 3. Funcoes_Inputs data rows (columns A/B) ------------------------------
$wsIn.Range("A2").Value = "calcular_eventos"
$wsIn.Range("B2").Value = "Pev_Tipico"

$wsIn.Range("A3").Value = "calcular_eventos"
$wsIn.Range("B3").Value = "Pev_Trajeto"

$wsIn.Range("A4").Value = "calcular_eventos"
$wsIn.Range("B4").Value = "Pev_DoenOcup"

$wsIn.Range("A5").Value = "calcular_eventos"
$wsIn.Range("B5").Value = "Pev_NRelac"

$wsIn.Range("A6").Value = "calcular_eventos"
$wsIn.Range("B6").Value = "Funcionarios"

# --- 4. Funcoes_Inputs column C (flags whether each input is external) -----
$wsIn.Range("C1").Value = "Param_Externo"
$wsIn.Range("C2").Formula = "=TRUE()"
$wsIn.Range("C3").Formula = "=TRUE()"
$wsIn.Range("C4").Formula = "=TRUE()"
$wsIn.Range("C5").Formula = "=TRUE()"
$wsIn.Range("C6").Formula = "=FALSE()"

# --- 5. Funcoes_Outputs data rows (columns A/B) -----------------------------
$wsOut.Range("A2").Value = "calcular_eventos"
$wsOut.Range("B2").Value = "Nev_Tipico"

$wsOut.Range("A3").Value = "calcular_eventos"
$wsOut.Range("B3").Value = "Nev_Trajeto"

$wsOut.Range("A4").Value = "calcular_eventos"
$wsOut.Range("B4").Value = "Nev_DoenOcup"

$wsOut.Range("A5").Value = "calcular_eventos"
$wsOut.Range("B5").Value = "Nev_NRelac"

# --- 6. Shade the first data row's "name" cell on each new sheet, matching
#        the look already used for similar catalogue tables elsewhere in the
#        workbook (e.g. Parametros!E2). -------------------------------------
$parametros.Range("E2").Copy()
$wsIn.Range("B2").PasteSpecial(-4122)   # xlPasteFormats
$wsOut.Range("B2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 7. Parametros: drop the obsolete trailing blank row (old row 18) ------
$parametros.Rows("18:18").Delete()

# --- 8. Refresh the AutoFilter so its range grows back to row 22 -----------
$parametros.AutoFilterMode = $false
$parametros.Range("A1:H22").AutoFilter()

# Keep the hidden _FilterDatabase defined name in sync with the new range.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Parametros!_FilterDatabase") {
        $n.RefersTo = "=Parametros!`$A`$1:`$H`$22"
    }
}

# --- 9. Restore Parametros as the selected/active sheet (unchanged from
#        before the edit). ---------------------------------------------------
$parametros.Activate()
$parametros.Range("A1").Select()
